$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.494.75"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.627.59"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'212.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.0621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'19.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "1.853.54"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "1.618.64"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "'63.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "26.480.78"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'214.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.86%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "'4.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'6.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +4.73%  "
$ws.Range("D25").Value = "'147.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.56%  "
$ws.Range("D29").Value = "'15.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "'2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "1.216.54"
$ws.Range("E36").Value = "  +5.63%  "
$ws.Range("E37").Value = "  +5.01%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'0.794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'2.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Value = "'0.793"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").Value = "1.762.43"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'92.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "'54.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").Value = "'0.0511"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "'7.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "'0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
